# The underlying data rows (2-7) were rotated: the last two records (old rows
# 6 and 7) were moved to the top (new rows 2 and 3), and the remaining
# records (old rows 2-5) shifted down by two rows (new rows 4-7).
#
# Concretely: new_row[N] = old_row[(N-2+4) mod 6 + 2]
#   old row 2 -> new row 4
#   old row 3 -> new row 5
#   old row 4 -> new row 6
#   old row 5 -> new row 7
#   old row 6 -> new row 2
#   old row 7 -> new row 3
#
# Rather than physically moving ranges (which risks disturbing formatting /
# helper cells), we just rewrite every data cell with its final value.
#
# The Startdatum/Slutdatum columns (Y/AA) hold plain text that looks like a
# date ("2019-09-24"). A bare `.Value = "2019-09-24"` assignment lets Excel
# auto-convert that into a real date serial number, which is not what the
# source file has (those cells are plain inline strings). Force them to stay
# text, then drop back to the default "Normal" style so no stray
# text-number-format style sticks around on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- New row 2 (was old row 6: id 80448314) ----
$ws.Range("A2").Value = 80448314
$ws.Range("B2").Value = 90665
$ws.Range("E2").Value = 4366
$ws.Range("F2").Value = "Skarp dropptaggsvamp"
$ws.Range("G2").Value = "Hydnellum peckii"
$ws.Range("H2").Value = "Banker"
$ws.Range("P2").Value = "En route, Dlr"
$ws.Range("Q2").Value = 501081.1366199313
$ws.Range("R2").Value = 6796554.212503371
$ws.Range("S2").Value = 5
Set-TextValue $ws.Range("Y2") "2019-09-24"
Set-TextValue $ws.Range("AA2") "2019-09-24"
$ws.Range("AW2").Value = "Andreas Öster"
$ws.Range("AX2").Value = "Andreas Öster"

# ---- New row 3 (was old row 7: id 90732145) ----
$ws.Range("A3").Value = 90732145
$ws.Range("B3").Value = 77605
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 967
$ws.Range("F3").Value = "Varglav"
$ws.Range("G3").Value = "Letharia vulpina"
$ws.Range("H3").Value = "(L.) Hue"
$ws.Range("P3").Value = "Jordikamäck, Dlr"
$ws.Range("Q3").Value = 501049.8414784005
$ws.Range("R3").Value = 6796606.094172257
$ws.Range("S3").Value = 5
Set-TextValue $ws.Range("Y3") "2020-10-28"
Set-TextValue $ws.Range("AA3") "2020-10-28"
$ws.Range("AW3").Value = "Andreas Öster"
$ws.Range("AX3").Value = "Andreas Öster"

# ---- New row 4 (was old row 2: id 73920018) ----
$ws.Range("A4").Value = 73920018
$ws.Range("B4").Value = 90653
$ws.Range("E4").Value = 4364
$ws.Range("F4").Value = "Dropptaggsvamp"
$ws.Range("G4").Value = "Hydnellum ferrugineum"
$ws.Range("H4").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q4").Value = 500981.982862098
$ws.Range("R4").Value = 6796631.058847403

# ---- New row 5 (was old row 3: id 73920027) ----
$ws.Range("A5").Value = 73920027
$ws.Range("B5").Value = 90653
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P5").Value = "Dalarna, Dlr"
$ws.Range("Q5").Value = 501185.0981308258
$ws.Range("R5").Value = 6796503.795656851
$ws.Range("S5").Value = 10
Set-TextValue $ws.Range("Y5") "2018-08-01"
Set-TextValue $ws.Range("AA5") "2018-08-01"
$ws.Range("AW5").Value = "Ville Pokela"
$ws.Range("AX5").Value = "Ville Pokela"

# ---- New row 6 (was old row 4: id 73920017) ----
$ws.Range("A6").Value = 73920017
$ws.Range("B6").Value = 88806
$ws.Range("E6").Value = 5685
$ws.Range("F6").Value = "Gullgröppa"
$ws.Range("G6").Value = "Pseudomerulius aureus"
$ws.Range("H6").Value = "(Fr.) Jülich"
$ws.Range("P6").Value = "Dalarna, Dlr"
$ws.Range("Q6").Value = 500982.9405652633
$ws.Range("R6").Value = 6796647.875423959
$ws.Range("S6").Value = 10
Set-TextValue $ws.Range("Y6") "2018-08-01"
Set-TextValue $ws.Range("AA6") "2018-08-01"
$ws.Range("AW6").Value = "Ville Pokela"
$ws.Range("AX6").Value = "Ville Pokela"

# ---- New row 7 (was old row 5: id 80081470) ----
$ws.Range("A7").Value = 80081470
$ws.Range("B7").Value = 88476
$ws.Range("E7").Value = 1962
$ws.Range("F7").Value = "Vaddporing"
$ws.Range("G7").Value = "Anomoporia kamtschatica"
$ws.Range("H7").Value = "(Parmasto) Bondartseva"
$ws.Range("P7").Value = "Mickelsjön, Dlr"
$ws.Range("Q7").Value = 500989.2359758026
$ws.Range("R7").Value = 6796507.578677795
Set-TextValue $ws.Range("Y7") "2019-09-24"
Set-TextValue $ws.Range("AA7") "2019-09-24"
